$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" column: copy the header formatting from the neighboring
# "sum" header (G1) onto H1, then set its text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# New data value for row 2 under the "Save" column.
$ws.Range("H2").Value = 0
